$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C1").Value = "ReportingManager"
$ws.Range("D1").Value = "ReportingManagerEmail"
$ws.Range("E1").Value = "DepartmentName"
$ws.Range("F1").Value = "ContactNo"
$ws.Range("G1").Value = "LocationOfficeName"
